# Update column F (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -1
    3  = 1
    4  = 2
    6  = 4
    7  = 7
    8  = 3
    9  = 1
    10 = -1
    11 = -2
    12 = -2
    13 = 4
    15 = -1
    17 = -2
    18 = -1
    19 = -6
    20 = 3
    22 = 3
    23 = -4
    24 = -3
    25 = -1
    26 = 4
    27 = 1
    28 = 2
    29 = -1
    30 = 2
    31 = 2
    32 = 6
    33 = -7
    34 = 2
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
